# Update measured results for worksheet "Planilha1" (fEntrada.xlsx)
# per commit: "Refeito os resultados com as medidas mais precisas e
# com a correcao da formula do fator de delaminacao pela area"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")
$ws.Activate()

# Recalculated Fd (B), Area_Del (C) and Fa (D) values for rows 2-22 (images d1..d21)
$ws.Cells.Item(2, 2).Value = 1.1056
$ws.Cells.Item(2, 3).Value = 0.1105
$ws.Cells.Item(2, 4).Value = 3.1318
$ws.Cells.Item(3, 2).Value = 1.1036
$ws.Cells.Item(3, 3).Value = 0.0827
$ws.Cells.Item(3, 4).Value = 2.3456
$ws.Cells.Item(4, 2).Value = 1.1161
$ws.Cells.Item(4, 3).Value = 0.0726
$ws.Cells.Item(4, 4).Value = 2.05706
$ws.Cells.Item(5, 2).Value = 1.0656
$ws.Cells.Item(5, 3).Value = 0.0828
$ws.Cells.Item(5, 4).Value = 2.3463
$ws.Cells.Item(6, 2).Value = 1.0817
$ws.Cells.Item(6, 3).Value = 0.0943
$ws.Cells.Item(6, 4).Value = 2.6733
$ws.Cells.Item(7, 2).Value = 1.0825
$ws.Cells.Item(7, 3).Value = 0.0919
$ws.Cells.Item(7, 4).Value = 2.0644
$ws.Cells.Item(8, 2).Value = 1.0518
$ws.Cells.Item(8, 3).Value = 0.0631
$ws.Cells.Item(8, 4).Value = 1.7904
$ws.Cells.Item(9, 2).Value = 1.0742
$ws.Cells.Item(9, 3).Value = 0.0675
$ws.Cells.Item(9, 4).Value = 1.9145
$ws.Cells.Item(10, 2).Value = 1.0688
$ws.Cells.Item(10, 3).Value = 0.0687
$ws.Cells.Item(10, 4).Value = 1.9476
$ws.Cells.Item(11, 2).Value = 1.0836
$ws.Cells.Item(11, 3).Value = 0.0688
$ws.Cells.Item(11, 4).Value = 1.9506
$ws.Cells.Item(12, 2).Value = 1.0817
$ws.Cells.Item(12, 3).Value = 0.0779
$ws.Cells.Item(12, 4).Value = 2.2097
$ws.Cells.Item(13, 2).Value = 1.09009
$ws.Cells.Item(13, 3).Value = 0.06405
$ws.Cells.Item(13, 4).Value = 1.8147
$ws.Cells.Item(14, 2).Value = 1.0783
$ws.Cells.Item(14, 3).Value = 0.0815
$ws.Cells.Item(14, 4).Value = 2.3113
$ws.Cells.Item(15, 2).Value = 1.0667
$ws.Cells.Item(15, 3).Value = 0.0832
$ws.Cells.Item(15, 4).Value = 2.3581
$ws.Cells.Item(16, 2).Value = 1.0658
$ws.Cells.Item(16, 3).Value = 0.0658
$ws.Cells.Item(16, 4).Value = 1.8652
$ws.Cells.Item(17, 2).Value = 1.0808
$ws.Cells.Item(17, 3).Value = 0.07901
$ws.Cells.Item(17, 4).Value = 2.2386
$ws.Cells.Item(18, 2).Value = 1.1157
$ws.Cells.Item(18, 3).Value = 0.1171
$ws.Cells.Item(18, 4).Value = 3.32001
$ws.Cells.Item(19, 2).Value = 1.0665
$ws.Cells.Item(19, 3).Value = 0.0708
$ws.Cells.Item(19, 4).Value = 2.0082
$ws.Cells.Item(20, 2).Value = 1.0651
$ws.Cells.Item(20, 3).Value = 0.0706
$ws.Cells.Item(20, 4).Value = 2.0028
$ws.Cells.Item(21, 2).Value = 1.0889
$ws.Cells.Item(21, 3).Value = 0.1053
$ws.Cells.Item(21, 4).Value = 2.9847
$ws.Cells.Item(22, 2).Value = 1.1478
$ws.Cells.Item(22, 3).Value = 0.1003
$ws.Cells.Item(22, 4).Value = 2.84306

# Update the view: zoom level and active selection
$ws.Range("E22").Select()
$excel.ActiveWindow.Zoom = 85
